# Generate Report for Handoff
#
# Updates the localization-status report after a fresh handoff run:
#   - Priority for the still-"low" rows (rows 4-7, the .md handoff/backhand
#     entries) is recomputed to "ht" on both the zh-cn and de-de sheets.
#   - Latest Handoff Datetime (column H) for those same rows is bumped to
#     the new handoff timestamp, per locale.
#   - The Overview sheet's "Latest HO Xliff Generate Date" column (G),
#     which mirrors the de-de handoff timestamp, is bumped to match.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$zhHandoffDate = "2016-09-07 06:45:27"
$deHandoffDate = "2016-09-07 06:45:32"

foreach ($r in 4..7) {
    # Priority column (E) recalculated from "low" to "ht"
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 5).Value = "ht"

    # Latest Handoff Datetime column (H), per locale
    $wsZhCn.Cells.Item($r, 8).Value = $zhHandoffDate
    $wsDeDe.Cells.Item($r, 8).Value = $deHandoffDate

    # Overview's Latest HO Xliff Generate Date column (G) tracks de-de
    $wsOverview.Cells.Item($r, 7).Value = $deHandoffDate
}

Write-Output "Report regenerated for handoff."
